# Update historicos/dolar_bevsa_uyu_temp.xlsx
# Insert two new rows of market data at the top of the historical table
# (rows 3 and 4), pushing older rows down, and drop the oldest row that
# falls off the bottom of the table (the former row 23 / 46031).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 3 (where data starts),
# shifting existing data rows 3-23 down to rows 5-25.
$ws.Range("A3:H4").EntireRow.Insert()

# New row 3: 2026-02-10 (serial 46063)
$ws.Cells.Item(3, 1).Value = 46063
$ws.Cells.Item(3, 2).Value = 38.469
$ws.Cells.Item(3, 3).Value = 38.469
$ws.Cells.Item(3, 4).Value = 38.48
$ws.Cells.Item(3, 5).Value = 38.45
$ws.Cells.Item(3, 6).Value = 38.49
$ws.Cells.Item(3, 7).Value = 28
$ws.Cells.Item(3, 8).Value = 14500000

# New row 4: 2026-02-09 (serial 46062)
$ws.Cells.Item(4, 1).Value = 46062
$ws.Cells.Item(4, 2).Value = 38.364
$ws.Cells.Item(4, 3).Value = 38.364
$ws.Cells.Item(4, 4).Value = 38.4
$ws.Cells.Item(4, 5).Value = 38.38
$ws.Cells.Item(4, 6).Value = 38.44
$ws.Cells.Item(4, 7).Value = 78
$ws.Cells.Item(4, 8).Value = 44000000

# Copy style (date format) from row 5 (shifted original row 3) to the new rows
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# The table now spans rows 3-25 (23 data rows); drop the oldest row that
# fell off the bottom (originally row 23, date 46031) so the table keeps
# the same number of rows as before plus the two newly added ones, minus
# the one that rolled off (net +1 row -> dimension A1:H24).
$ws.Range("A25:H25").EntireRow.Delete()
